$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 88.166664
$ws.Range("I11").Value = 88.166664
$ws.Range("K11").Value = 88.166664
$ws.Range("M11").Value = 51.833336
$ws.Range("H28").Value = 3468
$ws.Range("J28").Value = 4933.4614
$ws.Range("L28").Value = 4933.4614
$ws.Range("N28").Value = -5903.4614
$ws.Range("H40").Value = 6899.4
$ws.Range("J40").Value = 8749.25
$ws.Range("L40").Value = 8749.25
$ws.Range("N40").Value = -9099.25
$ws.Range("H74").Value = 11316.333
$ws.Range("I74").Value = 6975
$ws.Range("J74").Value = 19999
$ws.Range("K74").Value = 6975
$ws.Range("L74").Value = 19999
$ws.Range("M74").Value = -6039
$ws.Range("N74").Value = -21871
$ws.Range("H77").Value = 11316.333
$ws.Range("I77").Value = 6975
$ws.Range("J77").Value = 19999
$ws.Range("K77").Value = 34875
$ws.Range("L77").Value = 99995
$ws.Range("M77").Value = -30195
$ws.Range("N77").Value = -109355
$ws.Range("H87").Value = 97677
$ws.Range("J87").Value = 97677
$ws.Range("L87").Value = 97677
$ws.Range("N87").Value = -100173
$ws.Range("H90").Value = 97677
$ws.Range("J90").Value = 97677
$ws.Range("L90").Value = 293031
$ws.Range("N90").Value = -305511
$ws.Range("H98").Value = 580.8889
$ws.Range("I98").Value = 580.8889
$ws.Range("K98").Value = 580.8889
$ws.Range("M98").Value = 917.1111
$ws.Range("H100").Value = 3500.3333
$ws.Range("I100").Value = 3500.5
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 3500.5
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -2959.5
$ws.Range("N100").Value = -4582
$ws.Range("H116").Value = 3314.8333
$ws.Range("I116").Value = 2777.8
$ws.Range("K116").Value = 2777.8
$ws.Range("M116").Value = 664.1999999999998
$ws.Range("H122").Value = 580.8889
$ws.Range("I122").Value = 580.8889
$ws.Range("K122").Value = 1742.6667
$ws.Range("M122").Value = 707.3332999999998
$ws.Range("H125").Value = 2015
$ws.Range("I125").Value = 2015
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 18135
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -15675
$ws.Range("H141").Value = 2923.5715
$ws.Range("I141").Value = 2577.5
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 7732.5
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -2552.5
$ws.Range("N141").Value = -25360

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1822.8125
$ws.Range("I132").Value = 1744.3334
$ws.Range("K132").Value = 5233.0002
$ws.Range("M132").Value = -2703.0002
$ws.Range("H134").Value = 1810
$ws.Range("I134").Value = 1911.25
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 5733.75
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -3198.75
$ws.Range("N134").Value = -8070

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 3001
$ws.Range("I100").Value = 3001
$ws.Range("K100").Value = 9003
$ws.Range("M100").Value = -8192
$ws.Range("H131").Value = 2332.7856
$ws.Range("I131").Value = 1301.8
$ws.Range("K131").Value = 3905.4
$ws.Range("M131").Value = 1134.6

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3668.0667
$ws.Range("I102").Value = 3002.4443
$ws.Range("J102").Value = 4666.5
$ws.Range("K102").Value = 3002.4443
$ws.Range("L102").Value = 4666.5
$ws.Range("M102").Value = -1380.4443
$ws.Range("N102").Value = -7910.5
$ws.Range("H107").Value = 672.4
$ws.Range("I107").Value = 186.5
$ws.Range("J107").Value = 996.3333
$ws.Range("K107").Value = 186.5
$ws.Range("L107").Value = 996.3333
$ws.Range("M107").Value = 1733.5
$ws.Range("N107").Value = -4836.3333
$ws.Range("H113").Value = 9999.714
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 9999.714
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 9999.714
$ws.Range("N113").Value = -14339.714
$ws.Range("H132").Value = 2074.9285
$ws.Range("I132").Value = 1649.8572
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 4949.571599999999
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -2419.571599999999
$ws.Range("N132").Value = -12560

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1307.1428
$ws.Range("I7").Value = 1340
$ws.Range("J7").Value = 1225
$ws.Range("K7").Value = 1340
$ws.Range("L7").Value = 1225
$ws.Range("M7").Value = -1228
$ws.Range("N7").Value = -1449
$ws.Range("H20").Value = 502499.5
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 999999
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 999999
$ws.Range("M20").Value = -4774
$ws.Range("N20").Value = -1000451
$ws.Range("H122").Value = 3662.7334
$ws.Range("I122").Value = 3378.5
$ws.Range("J122").Value = 4799.6665
$ws.Range("K122").Value = 10135.5
$ws.Range("L122").Value = 14398.9995
$ws.Range("M122").Value = -7685.5
$ws.Range("N122").Value = -19298.9995
$ws.Range("H126").Value = 1307.1428
$ws.Range("I126").Value = 1340
$ws.Range("J126").Value = 1225
$ws.Range("K126").Value = 4020
$ws.Range("L126").Value = 3675
$ws.Range("M126").Value = -1550
$ws.Range("N126").Value = -8615
$ws.Range("H132").Value = 5983
$ws.Range("I132").Value = 6179.6
$ws.Range("K132").Value = 18538.8
$ws.Range("M132").Value = -16008.8

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3003.4
$ws.Range("I122").Value = 1740.5714
$ws.Range("K122").Value = 5221.7142
$ws.Range("M122").Value = -2771.7142
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0
$ws.Range("H126").Value = 4052.4707
$ws.Range("I126").Value = 1673.1111
$ws.Range("K126").Value = 5019.3333
$ws.Range("M126").Value = -2549.3333
